$d = $word.ActiveDocument

# The resume was cleaned up before building/publishing it: the two
# reviewer comments left on the "Add data filtering..." / "Build
# self-serve tools..." / "Automate unit tests..." bullet were removed.
# Deleting every comment also removes the commentRangeStart /
# commentRangeEnd / commentReference markers that anchored them, and
# drops the (now unused) comments part, which is why the hyperlink
# relationship ids shift down by one afterwards.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}
